$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the contents of rows 9-13 (columns A:AY) as follows:
#   new row 9  <- old row 11
#   new row 10 <- old row 12
#   new row 11 <- old row 10
#   new row 12 <- old row 13
#   new row 13 <- old row 9
# Row 14 keeps its own content (only column B changes).
#
# To do a safe in-place rotation we first stash each of the source rows
# (9-13) into scratch rows far below the used range, then copy each
# stashed row into its final destination, then clear the scratch rows.

$sourceRows = @(9, 10, 11, 12, 13)
$scratchBase = 2000

foreach ($r in $sourceRows) {
    $scratchRow = $scratchBase + $r
    $src = $ws.Range("A$r`:AY$r")
    $src.Copy()
    $dst = $ws.Range("A$scratchRow`:AY$scratchRow")
    $dst.PasteSpecial()
}

# mapping: destination row -> source row (stashed)
$mapping = @{
    9  = 11
    10 = 12
    11 = 10
    12 = 13
    13 = 9
}

# PasteSpecial on this engine skips blank source cells instead of
# overwriting the destination with blanks, so clear the destination rows
# first to avoid leaking stale values from the rows' original contents.
foreach ($destRow in @(9, 10, 11, 12, 13)) {
    $ws.Range("A$destRow`:AY$destRow").ClearContents()
}

foreach ($destRow in @(9, 10, 11, 12, 13)) {
    $srcRow = $mapping[$destRow]
    $scratchRow = $scratchBase + $srcRow
    $src = $ws.Range("A$scratchRow`:AY$scratchRow")
    $src.Copy()
    $dst = $ws.Range("A$destRow`:AY$destRow")
    $dst.PasteSpecial()
}

# clear the scratch rows
foreach ($r in $sourceRows) {
    $scratchRow = $scratchBase + $r
    $ws.Range("A$scratchRow`:AY$scratchRow").ClearContents()
}

$excel.CutCopyMode = 0

# Column B ("Taxonsorteringsordning") gets fresh values independent of the
# row rotation above.
$ws.Range("B9").Value = 90814
$ws.Range("B10").Value = 103742
$ws.Range("B11").Value = 90806
$ws.Range("B12").Value = 103781
$ws.Range("B13").Value = 90806
$ws.Range("B14").Value = 90817
